$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.921897292137146
$ws.Range("B1").Value = 1.629784226417542
$ws.Range("D1").Value = 1.636160254478455
$ws.Range("E1").Value = 1.068482756614685
